# Auto-generated Excel COM-interop script
# Applies numeric value updates to the 'Hyperion_Profits' style Leve tables
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR per the target diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 922
$ws.Range("I11").Value = 922
$ws.Range("K11").Value = 922
$ws.Range("M11").Value = -782
$ws.Range("H29").Value = 4334.3335
$ws.Range("J29").Value = 3000
$ws.Range("L29").Value = 9000
$ws.Range("N29").Value = -9562

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8968.723
$ws.Range("J32").Value = 16748.76
$ws.Range("L32").Value = 16748.76
$ws.Range("N32").Value = -17322.76
$ws.Range("H45").Value = 5141188.5
$ws.Range("I45").Value = 7993835.5
$ws.Range("J45").Value = 6424.6
$ws.Range("K45").Value = 7993835.5
$ws.Range("L45").Value = 6424.6
$ws.Range("M45").Value = -7993458.5
$ws.Range("N45").Value = -7178.6
$ws.Range("H61").Value = 2644.5
$ws.Range("I61").Value = 2308.3572
$ws.Range("K61").Value = 2308.3572
$ws.Range("M61").Value = -2096.3572
$ws.Range("H132").Value = 2329.7646
$ws.Range("I132").Value = 2225.375
$ws.Range("K132").Value = 6676.125
$ws.Range("M132").Value = -4146.125
$ws.Range("H136").Value = 2644.5
$ws.Range("I136").Value = 2308.3572
$ws.Range("K136").Value = 6925.071599999999
$ws.Range("M136").Value = -4375.071599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 302
$ws.Range("I12").Value = 302
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 302
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -134
$ws.Range("N12").ClearContents()
$ws.Range("H14").Value = 2166.5
$ws.Range("I14").Value = 1000
$ws.Range("J14").Value = 3333
$ws.Range("K14").Value = 1000
$ws.Range("L14").Value = 3333
$ws.Range("M14").Value = -828
$ws.Range("N14").Value = -3677
$ws.Range("H94").Value = 2130794.5
$ws.Range("I94").Value = 2942411.5
$ws.Range("J94").Value = 8104.077
$ws.Range("K94").Value = 2942411.5
$ws.Range("L94").Value = 8104.077
$ws.Range("M94").Value = -2941960.5
$ws.Range("N94").Value = -9006.077000000001
$ws.Range("H99").Value = 6804986
$ws.Range("I99").Value = 11906351
$ws.Range("J99").Value = 3165.5557
$ws.Range("K99").Value = 11906351
$ws.Range("L99").Value = 3165.5557
$ws.Range("M99").Value = -11904853
$ws.Range("N99").Value = -6161.5557

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 16246.75
$ws.Range("I58").Value = 17854.143
$ws.Range("K58").Value = 17854.143
$ws.Range("M58").Value = -17651.143
$ws.Range("H136").Value = 16246.75
$ws.Range("I136").Value = 17854.143
$ws.Range("K136").Value = 53562.429
$ws.Range("M136").Value = -51012.429
$ws.Range("H137").Value = 50526.668
$ws.Range("J137").Value = 50526.668
$ws.Range("L137").Value = 50526.668
$ws.Range("N137").Value = -60726.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1461.1364
$ws.Range("I34").Value = 217.22223
$ws.Range("J34").Value = 2322.3076
$ws.Range("K34").Value = 651.66669
$ws.Range("L34").Value = 6966.9228
$ws.Range("M34").Value = -567.66669
$ws.Range("N34").Value = -7134.9228
$ws.Range("H39").Value = 2188.8
$ws.Range("J39").Value = 2286.25
$ws.Range("L39").Value = 6858.75
$ws.Range("N39").Value = -7446.75
$ws.Range("H55").Value = 35903.633
$ws.Range("I55").Value = 338.875
$ws.Range("K55").Value = 1016.625
$ws.Range("M55").Value = -839.625
$ws.Range("H137").Value = 3674.5
$ws.Range("I137").Value = 3800
$ws.Range("J137").Value = 3599.2
$ws.Range("K137").Value = 11400
$ws.Range("L137").Value = 10797.6
$ws.Range("M137").Value = -6300
$ws.Range("N137").Value = -20997.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("H10").Value = 2003.5
$ws.Range("J10").Value = 2003.5
$ws.Range("L10").Value = 2003.5
$ws.Range("N10").Value = -2341.5
$ws.Range("H11").Value = 1376.6
$ws.Range("I11").Value = 1550
$ws.Range("J11").Value = 1261
$ws.Range("K11").Value = 1550
$ws.Range("L11").Value = 1261
$ws.Range("M11").Value = -1411
$ws.Range("N11").Value = -1539
$ws.Range("H12").Value = 4999.3335
$ws.Range("H13").Value = 1098
$ws.Range("J13").Value = 1162.3334
$ws.Range("L13").Value = 1162.3334
$ws.Range("N13").Value = -1440.3334
$ws.Range("H122").Value = 308931.53
$ws.Range("I122").Value = 406439.7
$ws.Range("K122").Value = 1219319.1
$ws.Range("M122").Value = -1216869.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1495.1111
$ws.Range("J16").Value = 2895.2856
$ws.Range("L16").Value = 2895.2856
$ws.Range("N16").Value = -3235.2856
$ws.Range("H68").Value = 999
$ws.Range("J68").Value = 999
$ws.Range("L68").Value = 999
$ws.Range("N68").Value = -2497
$ws.Range("H71").Value = 999
$ws.Range("J71").Value = 999
$ws.Range("L71").Value = 4995
$ws.Range("N71").Value = -12483
$ws.Range("H93").Value = 11116099
$ws.Range("I93").Value = 12822084
$ws.Range("J93").Value = 27199.75
$ws.Range("K93").Value = 12822084
$ws.Range("L93").Value = 27199.75
$ws.Range("M93").Value = -12820836
$ws.Range("N93").Value = -29695.75
$ws.Range("H100").Value = 3085.625
$ws.Range("I100").Value = 2947.75
$ws.Range("K100").Value = 2947.75
$ws.Range("M100").Value = -2406.75
$ws.Range("H136").Value = 76667.89
$ws.Range("I136").Value = 157538.69
$ws.Range("K136").Value = 472616.07
$ws.Range("M136").Value = -470066.07

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 188149.4
$ws.Range("J4").Value = 1915.6666
$ws.Range("L4").Value = 1915.6666
$ws.Range("N4").Value = -2141.6666
$ws.Range("H8").Value = 5000000
$ws.Range("I8").Value = 5000000
$ws.Range("K8").Value = 5000000
$ws.Range("M8").Value = -4999860
$ws.Range("H33").Value = 13440
$ws.Range("I33").Value = 6637.5
$ws.Range("K33").Value = 6637.5
$ws.Range("M33").Value = -6387.5
$ws.Range("H36").Value = 13440
$ws.Range("I36").Value = 6637.5
$ws.Range("K36").Value = 6637.5
$ws.Range("M36").Value = -6387.5
$ws.Range("H62").Value = 7309
$ws.Range("I62").Value = 1366.5
$ws.Range("J62").Value = 8272.647999999999
$ws.Range("K62").Value = 1366.5
$ws.Range("L62").Value = 8272.647999999999
$ws.Range("M62").Value = -742.5
$ws.Range("N62").Value = -9520.647999999999
$ws.Range("H65").Value = 7309
$ws.Range("I65").Value = 1366.5
$ws.Range("J65").Value = 8272.647999999999
$ws.Range("K65").Value = 6832.5
$ws.Range("L65").Value = 41363.24
$ws.Range("M65").Value = -3712.5
$ws.Range("N65").Value = -47603.24
$ws.Range("H103").Value = 44951.25
$ws.Range("J103").Value = 44951.25
$ws.Range("L103").Value = 44951.25
$ws.Range("N103").Value = -47295.25
$ws.Range("H107").Value = 31250948
$ws.Range("I107").Value = 40000344
$ws.Range("K107").Value = 120001032
$ws.Range("M107").Value = -119999112
$ws.Range("H113").Value = 876.0417
$ws.Range("I113").Value = 543.4545000000001
$ws.Range("J113").Value = 1157.4615
$ws.Range("K113").Value = 1630.3635
$ws.Range("L113").Value = 3472.3845
$ws.Range("M113").Value = 539.6364999999998
$ws.Range("N113").Value = -7812.3845
$ws.Range("H132").Value = 24070632
$ws.Range("I132").Value = 26317780
$ws.Range("J132").Value = 2722721.5
$ws.Range("K132").Value = 78953340
$ws.Range("L132").Value = 8168164.5
$ws.Range("M132").Value = -78950810
$ws.Range("N132").Value = -8173224.5
$ws.Range("H136").Value = 2184.2856
$ws.Range("I136").Value = 1531.6666
$ws.Range("K136").Value = 4594.9998
$ws.Range("M136").Value = -2044.9998
